# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (keeps its sheetPr/header/column styling)
#    and place the copy immediately before "总计", then rename it "2022-Q1".
# 2. Overwrite the fund-holding rows with the 2022-Q1 figures (21 funds).
# 3. Insert a new "2022-Q1" row at the top of the "总计" summary sheet and
#    renumber the index column for the rows that shift down.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet, positioned right before "总计"
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($wb.Worksheets.Item("总计"))
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# NOTE: any worksheet handle obtained *before* the Copy() above (which
# inserts a brand-new sheet into the collection) must be re-resolved by
# name afterwards - a handle grabbed earlier tracks a *position* that the
# insertion shifts, so it would silently now point at the new sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# The template has 20 fund rows (rows 2-21); 2022-Q1 needs 21 (rows 2-22),
# so clone the styling of the last data row down into row 22 first.
$q1.Range("A21:H21").Copy()
$q1.Range("A22:H22").PasteSpecial(-4122)

# Force text storage (matches the source data: codes / ratios are strings,
# not numbers) without disturbing the existing header/index styling.
$q1.Range("B2:B22").NumberFormat = "@"
$q1.Range("D2:G22").NumberFormat = "@"

$fundData = @(
    @(0, "160527", "博时研究优选3年封闭运作灵活配置混合A", "20.12", "95.51", "6.18", "1.2434", 7),
    @(1, "011855", "银华长荣混合型证券投资基金", "11.96", "67.01", "4.83", "0.5777", 7),
    @(2, "000264", "博时内需增长混合", "3.87", "75.26", "4.22", "0.1633", 9),
    @(3, "014107", "博时品质生活混合A", "4.11", "69.12", "3.67", "0.1508", 7),
    @(4, "008978", "银华长丰混合", "2.51", "70.68", "4.99", "0.1252", 6),
    @(5, "005265", "博时厚泽回报灵活配置混合A", "2.19", "91.85", "3.66", "0.0802", 4),
    @(6, "012153", "博时研究慧选混合型证券投资基金A", "1.63", "75.28", "4.49", "0.0732", 3),
    @(7, "160528", "博时研究优选3年封闭运作灵活配置混合C", "0.69", "95.51", "6.18", "0.0426", 7),
    @(8, "562500", "华夏中证机器人ETF", "1.35", "99.22", "2.27", "0.0306", 10),
    @(9, "005266", "博时厚泽回报灵活配置混合C", "0.64", "91.85", "3.66", "0.0234", 4),
    @(10, "004677", "博时战略新兴产业混合", "0.41", "89.27", "5.24", "0.0215", 5),
    @(11, "009700", "长江添利混合A", "1.90", "20.46", "1.13", "0.0215", 4),
    @(12, "159770", "天弘中证机器人ETF", "0.87", "99.47", "2.28", "0.0198", 10),
    @(13, "009701", "长江添利混合C", "1.40", "20.46", "1.13", "0.0158", 4),
    @(14, "562360", "银华中证机器人ETF", "0.68", "96.94", "2.23", "0.0152", 10),
    @(15, "010663", "长江均衡成长混合A", "0.26", "85.90", "5.64", "0.0147", 1),
    @(16, "003659", "山西证券策略精选灵活配置混合", "0.31", "84.52", "4.75", "0.0147", 1),
    @(17, "012154", "博时研究慧选混合型证券投资基金C", "0.21", "75.28", "4.49", "0.0094", 3),
    @(18, "014108", "博时品质生活混合C", "0.14", "69.12", "3.67", "0.0051", 7),
    @(19, "010664", "长江均衡成长混合C", "0.05", "85.90", "5.64", "0.0028", 1),
    @(20, "004696", "东兴量化优享灵活配置混合", "0.03", "68.71", "3.10", "0.0009", 5)
)

$r = 2
foreach ($row in $fundData) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2) Insert the 2022-Q1 summary row at the top of "总计" and renumber
# ------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 21
$totalSheet.Cells.Item(2, 4).Value = 2.65

for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
